$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.1818801724491279
$ws.Range("J2").Value = 0.2500781318045117
$ws.Range("M2").Value = 0.1798956666666667
$ws.Range("N2").Value = 0.539687
$ws.Range("O2").Value = 0.01070918551864568
$ws.Range("P2").Value = 0.01088280728907136
$ws.Range("Q2").Value = 0.005853924923777778
$ws.Range("R2").Value = 0.05268532431400001
$ws.Range("S2").Value = 0.001947788508920979
$ws.Range("T2").Value = 0.002721552115639489
$ws.Range("I3").Value = 0.1818801724491279
$ws.Range("J3").Value = 0.2500781318045117
$ws.Range("N3").Value = 47.39813
$ws.Range("O3").Value = 0.9405365839956962
$ws.Range("P3").Value = 0.9557849543390003
$ws.Range("Q3").Value = 0.5141222496511112
$ws.Range("S3").Value = 0.1710649560918509
$ws.Range("T3").Value = 0.2390209157879577
$ws.Range("I4").Value = 0.1818801724491279
$ws.Range("J4").Value = 0.2500781318045117
$ws.Range("K4").Value = 2
$ws.Range("L4").Value = 0.6666666666666666
$ws.Range("M4").Value = 0.01499966666666667
$ws.Range("N4").Value = 0.044999
$ws.Range("O4").Value = 0.0008929298633347419
$ws.Range("P4").Value = 0.0009074064137192897
$ws.Range("Q4").Value = 0.0004880991531111111
$ws.Range("R4").Value = 0.004392892378
$ws.Range("S4").Value = 0.0001624062375282991
$ws.Range("T4").Value = 0.0002269225007303517
$ws.Range("I5").Value = 0.1818801724491279
$ws.Range("J5").Value = 0.2500781318045117
$ws.Range("M5").Value = 0.8039865
$ws.Range("N5").Value = 1.607973
$ws.Range("O5").Value = 0.04786130062232345
$ws.Range("P5").Value = 0.03242483195820901
$ws.Range("Q5").Value = 0.026162256701
$ws.Range("R5").Value = 0.156973540206
$ws.Range("S5").Value = 0.008705021610827741
$ws.Range("T5").Value = 0.008108741400184135
$ws.Range("E6").Value = 1
$ws.Range("F6").Value = 0.5
$ws.Range("G6").Value = 0.146372
$ws.Range("H6").Value = 0.292744
$ws.Range("I6").Value = 0.8181198275508721
$ws.Range("J6").Value = 0.7499218681954883
$ws.Range("M6").Value = 0.1798956666666667
$ws.Range("N6").Value = 0.539687
$ws.Range("O6").Value = 0.01070918551864568
$ws.Range("P6").Value = 0.01088280728907136
$ws.Range("Q6").Value = 0.02633168852133334
$ws.Range("R6").Value = 0.157990131128
$ws.Range("S6").Value = 0.008761397009724699
$ws.Range("T6").Value = 0.008161255173431876
$ws.Range("E7").Value = 1
$ws.Range("F7").Value = 0.5
$ws.Range("G7").Value = 0.146372
$ws.Range("H7").Value = 0.292744
$ws.Range("I7").Value = 0.8181198275508721
$ws.Range("J7").Value = 0.7499218681954883
$ws.Range("N7").Value = 47.39813
$ws.Range("O7").Value = 0.9405365839956962
$ws.Range("P7").Value = 0.9557849543390003
$ws.Range("Q7").Value = 2.312586361453334
$ws.Range("R7").Value = 13.87551816872
$ws.Range("S7").Value = 0.7694716279038453
$ws.Range("T7").Value = 0.7167640385510426
$ws.Range("E8").Value = 1
$ws.Range("F8").Value = 0.5
$ws.Range("G8").Value = 0.146372
$ws.Range("H8").Value = 0.292744
$ws.Range("I8").Value = 0.8181198275508721
$ws.Range("J8").Value = 0.7499218681954883
$ws.Range("K8").Value = 2
$ws.Range("L8").Value = 0.6666666666666666
$ws.Range("M8").Value = 0.01499966666666667
$ws.Range("N8").Value = 0.044999
$ws.Range("O8").Value = 0.0008929298633347419
$ws.Range("P8").Value = 0.0009074064137192897
$ws.Range("Q8").Value = 0.002195531209333333
$ws.Range("R8").Value = 0.013173187256
$ws.Range("S8").Value = 0.0007305236258064429
$ws.Range("T8").Value = 0.0006804839129889379
$ws.Range("E9").Value = 1
$ws.Range("F9").Value = 0.5
$ws.Range("G9").Value = 0.146372
$ws.Range("H9").Value = 0.292744
$ws.Range("I9").Value = 0.8181198275508721
$ws.Range("J9").Value = 0.7499218681954883
$ws.Range("M9").Value = 0.8039865
$ws.Range("N9").Value = 1.607973
$ws.Range("O9").Value = 0.04786130062232345
$ws.Range("P9").Value = 0.03242483195820901
$ws.Range("Q9").Value = 0.117681111978
$ws.Range("R9").Value = 0.470724447912
$ws.Range("S9").Value = 0.0391562790114957
$ws.Range("T9").Value = 0.2544943929043076
